# Se corrigió item en SelectCampania Modif de DataSourceMotor
# Se corrigió error con Item Refacturacion

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before column F (old F.. shifts right by 2, e.g. old F -> H)
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels for the two newly inserted columns
$ws.Range("F1").Value = "CodigoAgente"
$ws.Range("G1").Value = "NUM_GRUPO"

# New data values for row 2 in the inserted columns
$ws.Range("F2").Value = 2302
$ws.Range("G2").Value = "Mattioli"

# Center-align the new data cells (matches the added style with centered alignment)
$dataRange = $ws.Range("F2:G2")
$dataRange.Font.Color = 0
$dataRange.HorizontalAlignment = -4108  # xlCenter

# Corrected account/reference number in E2 (Item Refacturacion fix)
$ws.Range("E2").Value = 2617100594

# Update active selection to the newly added cells
$ws.Range("F2:G2").Select() | Out-Null
